$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 184, shifting existing rows 184:200 down to 185:201
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new weekly data point
$ws.Cells.Item(184, 1).Value = 11
$ws.Cells.Item(184, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(184, 3).Value = "Bíobío"
$ws.Cells.Item(184, 4).Value = 45223
$ws.Cells.Item(184, 5).Value = 8
$ws.Cells.Item(184, 6).Value = "Fruta"
$ws.Cells.Item(184, 7).Value = 100108
$ws.Cells.Item(184, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(184, 9).Value = 100108002
$ws.Cells.Item(184, 10).Value = "Mango"
$ws.Cells.Item(184, 11).Value = "Sin especificar"
$ws.Cells.Item(184, 12).Value = "Primera"
$ws.Cells.Item(184, 13).Value = 200
$ws.Cells.Item(184, 14).Value = 9500
$ws.Cells.Item(184, 15).Value = 10000
$ws.Cells.Item(184, 16).Value = 9750
$ws.Cells.Item(184, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(184, 18).Value = "Brasil"
$ws.Cells.Item(184, 19).Value = 2438
$ws.Cells.Item(184, 20).Value = 4
